$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 11.14494766666667
$ws.Range("H2").Value = 33.434843
$ws.Range("I2").Value = 0.1279818847384872
$ws.Range("J2").Value = 0.1279818847384872
$ws.Range("O2").Value = 0.858667536176972
$ws.Range("P2").Value = 0.858667536176972
$ws.Range("Q2").Value = 1.751172192020333
$ws.Range("R2").Value = 15.760549728183
$ws.Range("S2").Value = 0.109893889643682
$ws.Range("T2").Value = 0.109893889643682

# Row 3
$ws.Range("G3").Value = 11.14494766666667
$ws.Range("H3").Value = 33.434843
$ws.Range("I3").Value = 0.1279818847384872
$ws.Range("J3").Value = 0.1279818847384872
$ws.Range("Q3").Value = 0.2882343515378889
$ws.Range("R3").Value = 2.594109163841
$ws.Range("S3").Value = 0.01808799509480518
$ws.Range("T3").Value = 0.01808799509480518

# Row 4
$ws.Range("I4").Value = 0.5307607770439682
$ws.Range("J4").Value = 0.5307607770439681
$ws.Range("O4").Value = 0.858667536176972
$ws.Range("P4").Value = 0.858667536176972
$ws.Range("Q4").Value = 7.262383385536999
$ws.Range("S4").Value = 0.4557470487237193
$ws.Range("T4").Value = 0.4557470487237192

# Row 5
$ws.Range("I5").Value = 0.5307607770439682
$ws.Range("J5").Value = 0.5307607770439681
$ws.Range("S5").Value = 0.07501372832024883
$ws.Range("T5").Value = 0.07501372832024882

# Row 6
$ws.Range("G6").Value = 29.71744933333333
$ws.Range("H6").Value = 89.152348
$ws.Range("I6").Value = 0.3412573382175446
$ws.Range("J6").Value = 0.3412573382175446
$ws.Range("O6").Value = 0.858667536176972
$ws.Range("P6").Value = 0.858667536176972
$ws.Range("Q6").Value = 4.669413661398666
$ws.Range("R6").Value = 42.024722952588
$ws.Range("S6").Value = 0.2930265978095707
$ws.Range("T6").Value = 0.2930265978095707

# Row 7
$ws.Range("G7").Value = 29.71744933333333
$ws.Range("H7").Value = 89.152348
$ws.Range("I7").Value = 0.3412573382175446
$ws.Range("J7").Value = 0.3412573382175446
$ws.Range("Q7").Value = 0.7685625804751112
$ws.Range("R7").Value = 6.917063224276
$ws.Range("S7").Value = 0.04823074040797393
$ws.Range("T7").Value = 0.04823074040797393
